$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price (column D) value would be auto-parsed as a number by
# Excel's input coercion (e.g. "0.5311"). For those we briefly mark the cell
# as Text ("@") before assigning, then restore the "Normal" style so the
# saved file keeps the original (unstyled) cell formatting.
$textCells = @(
    @{Addr='D6'; Value='0.5311'}
    @{Addr='D8'; Value='0.2622'}
    @{Addr='D9'; Value='0.06330'}
    @{Addr='D11'; Value='0.07805'}
    @{Addr='D12'; Value='4.517'}
    @{Addr='D15'; Value='0.5488'}
    @{Addr='D17'; Value='65.36'}
    @{Addr='D20'; Value='4.588'}
    @{Addr='D21'; Value='190.95'}
    @{Addr='D23'; Value='6.009'}
    @{Addr='D24'; Value='1.007'}
    @{Addr='D26'; Value='0.1219'}
    @{Addr='D27'; Value='7.187'}
    @{Addr='D29'; Value='1.473'}
    @{Addr='D30'; Value='0.05737'}
    @{Addr='D32'; Value='3.548'}
    @{Addr='D33'; Value='3.261'}
    @{Addr='D34'; Value='1.587'}
    @{Addr='D36'; Value='2.422'}
    @{Addr='D38'; Value='0.5743'}
    @{Addr='D44'; Value='103.72'}
    @{Addr='D46'; Value='56.82'}
    @{Addr='D49'; Value='0.4359'}
    @{Addr='D50'; Value='7.861'}
    @{Addr='D51'; Value='0.05152'}
)
foreach ($item in $textCells) {
    $cell = $ws.Range($item.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.Style = "Normal"
}

# Cells whose new Price value already round-trips as text unambiguously
# (contains more than one ".", or other non-numeric characters).
$plainCells = @(
    @{Addr='D2'; Value='26.174.27'}
    @{Addr='D3'; Value='1.653.58'}
    @{Addr='D13'; Value='1.624.55'}
    @{Addr='D14'; Value='1.881.36'}
    @{Addr='D16'; Value='0.0₅8161'}
    @{Addr='D18'; Value='26.142.14'}
    @{Addr='D43'; Value='1.039.84'}
    @{Addr='D45'; Value='1.793.90'}
)
foreach ($item in $plainCells) {
    $ws.Range($item.Addr).Value = $item.Value
}

# Volume(1h) column (E) updates - always text (percent sign + padding spaces).
$volumeCells = @(
    @{Addr='E2'; Value='  +0.51%  '}
    @{Addr='E3'; Value='  +0.15%  '}
    @{Addr='E6'; Value='  +0.48%  '}
    @{Addr='E7'; Value='  +0.41%  '}
    @{Addr='E8'; Value='  +0.09%  '}
    @{Addr='E9'; Value='  +0.39%  '}
    @{Addr='E10'; Value='  -0.27%  '}
    @{Addr='E11'; Value='  +0.78%  '}
    @{Addr='E12'; Value='  +0.75%  '}
    @{Addr='E13'; Value='  -0.58%  '}
    @{Addr='E14'; Value='  +0.24%  '}
    @{Addr='E15'; Value='  +0.42%  '}
    @{Addr='E16'; Value='  +0.73%  '}
    @{Addr='E17'; Value='  +0.43%  '}
    @{Addr='E18'; Value='  +0.31%  '}
    @{Addr='E19'; Value='  +0.43%  '}
    @{Addr='E20'; Value='  +0.45%  '}
    @{Addr='E21'; Value='  -1.04%  '}
    @{Addr='E22'; Value='  +0.18%  '}
    @{Addr='E23'; Value='  +0.20%  '}
    @{Addr='E24'; Value='  +0.44%  '}
    @{Addr='E25'; Value='  +4.02%  '}
    @{Addr='E26'; Value='  -2.03%  '}
    @{Addr='E27'; Value='  -1.24%  '}
    @{Addr='E28'; Value='  -1.76%  '}
    @{Addr='E29'; Value='  +4.38%  '}
    @{Addr='E30'; Value='  -3.56%  '}
    @{Addr='E32'; Value='  +1.36%  '}
    @{Addr='E33'; Value='  +0.38%  '}
    @{Addr='E34'; Value='  +3.02%  '}
    @{Addr='E35'; Value='  +1.83%  '}
    @{Addr='E36'; Value='  +0.37%  '}
    @{Addr='E37'; Value='  +0.39%  '}
    @{Addr='E38'; Value='  +1.44%  '}
    @{Addr='E39'; Value='  -0.63%  '}
    @{Addr='E40'; Value='  -1.53%  '}
    @{Addr='E41'; Value='  +0.09%  '}
    @{Addr='E42'; Value='  +0.48%  '}
    @{Addr='E43'; Value='  +3.21%  '}
    @{Addr='E44'; Value='  +2.76%  '}
    @{Addr='E45'; Value='  +0.07%  '}
    @{Addr='E46'; Value='  -0.07%  '}
    @{Addr='E47'; Value='  -1.77%  '}
    @{Addr='E48'; Value='  -0.24%  '}
    @{Addr='E49'; Value='  +1.73%  '}
    @{Addr='E50'; Value='  -0.17%  '}
    @{Addr='E51'; Value='  +0.03%  '}
)
foreach ($item in $volumeCells) {
    $ws.Range($item.Addr).Value = $item.Value
}
